$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 362
$wsExhibition.Range("F4").Value = 2967
$wsExhibition.Range("F6").Value = 619

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 362
$wsAll.Range("F6").Value = 2967
$wsAll.Range("F8").Value = 619
